$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Google")
$ws2 = $wb.Worksheets.Item("DuckDuckGo")

$ws1.Activate()
$ws1.Range("D7").Select()

$ws2.Activate()
$ws2.Range("A4").Value = "ruby is the best programming language"
$ws2.Columns.Item(1).ColumnWidth = 36
$ws2.Range("C9").Select()
